# Updates cryptos list figures (price/volume columns) for the sheet, and
# swaps the Dai / LEO rows (26 and 27) to reflect the new ranking order.
#
# Note: several "Price" values are plain decimal numbers (e.g. "580.74").
# Because column D holds these as text (to match the "NN.NNN.NN"-style
# thousands grouping used elsewhere in the column), we force a temporary
# text number format before assigning such values so Excel does not
# silently convert them to numeric cells, then restore the default style
# so the cell formatting matches the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.643.12'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '3.447.14'
$ws.Range("E3").Value = '  +2.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.10%  '
$ws.Range("D7").Value = '3.448.73'
$ws.Range("E7").Value = '  +2.40%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.78'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("E11").Value = '  +3.55%  '
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").Value = '4.036.83'
$ws.Range("E13").Value = '  +2.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.64%  '
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000175'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.37%  '
$ws.Range("D17").Value = '3.450.64'
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("D18").Value = '61.758.09'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.27'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = '  +1.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.566'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.63%  '
$ws.Range("D24").Value = '3.589.03'
$ws.Range("E24").Value = '  +2.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.88%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("E28").Value = '  +0.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.182'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.73'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.68%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  -13.31%  '
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("E34").Value = '  +1.40%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("E38").Value = '  +0.57%  '
$ws.Range("E39").Value = '  +1.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0786'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.47%  '
$ws.Range("E43").Value = '  +2.25%  '
$ws.Range("E44").Value = '  +2.45%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.96%  '
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("D48").Value = '2.606.76'
$ws.Range("E48").Value = '  +6.18%  '
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.76%  '
